# Weekly update: insert a new price record for "Poroto granado" (Femacal de
# La Calera) ahead of the existing rows, pushing the rest of the table down
# by one row (dimension grows from A1:R183 to A1:R184).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 86; Excel shifts rows 86:183 down to 87:184.
$ws.Rows.Item(86).Insert()

# Populate the newly inserted row 86 with the new record.
$ws.Cells.Item(86, 1).Value = 3
$ws.Cells.Item(86, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(86, 3).Value = "Coquimbo"
$ws.Cells.Item(86, 4).Value = 44679
$ws.Cells.Item(86, 5).Value = 5
$ws.Cells.Item(86, 6).Value = 100112030
$ws.Cells.Item(86, 7).Value = "Poroto granado"
$ws.Cells.Item(86, 8).Value = "Sin especificar"
$ws.Cells.Item(86, 9).Value = "Primera"
$ws.Cells.Item(86, 10).Value = 73
$ws.Cells.Item(86, 11).Value = 23000
$ws.Cells.Item(86, 12).Value = 24000
$ws.Cells.Item(86, 13).Value = 23521
$ws.Cells.Item(86, 14).Value = "$/malla 25 kilos"
$ws.Cells.Item(86, 15).Value = "Provincia de Talca"
$ws.Cells.Item(86, 16).Value = 941
$ws.Cells.Item(86, 17).Value = 25
$ws.Cells.Item(86, 18).Value = "Hortaliza"
